# Auto-generated Excel COM-interop edit script
# Applies per-cell numeric updates to the Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# as produced by the scheduled market-data refresh runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("H18").Value2 = 9000
$ws.Range("I18").Value2 = 8333.333000000001
$ws.Range("K18").Value2 = 8333.333000000001
$ws.Range("M18").Value2 = -8049.333000000001

$ws = $wb.Worksheets.Item(1)
$ws.Range("H28").Value2 = 967.625
$ws.Range("I28").Value2 = 689.4
$ws.Range("K28").Value2 = 689.4
$ws.Range("M28").Value2 = -204.4

$ws = $wb.Worksheets.Item(1)
$ws.Range("H33").Value2 = 2531.9312
$ws.Range("I33").Value2 = 2642.4092
$ws.Range("J33").Value2 = 2184.7144
$ws.Range("K33").Value2 = 2642.4092
$ws.Range("L33").Value2 = 2184.7144
$ws.Range("M33").Value2 = -2413.4092
$ws.Range("N33").Value2 = -2642.7144

$ws = $wb.Worksheets.Item(1)
$ws.Range("H43").Value2 = 3535.5715
$ws.Range("I43").Value2 = 1200.1
$ws.Range("J43").Value2 = 9374.25
$ws.Range("K43").Value2 = 1200.1
$ws.Range("L43").Value2 = 9374.25
$ws.Range("M43").Value2 = -1131.1
$ws.Range("N43").Value2 = -9512.25

$ws = $wb.Worksheets.Item(1)
$ws.Range("H62").Value2 = 4619.8945
$ws.Range("I62").Value2 = 2816.7273
$ws.Range("K62").Value2 = 2816.7273
$ws.Range("M62").Value2 = -2192.7273

$ws = $wb.Worksheets.Item(1)
$ws.Range("H65").Value2 = 4619.8945
$ws.Range("I65").Value2 = 2816.7273
$ws.Range("K65").Value2 = 14083.6365
$ws.Range("M65").Value2 = -10963.6365

$ws = $wb.Worksheets.Item(1)
$ws.Range("H135").Value2 = 7985.5
$ws.Range("I135").Value2 = 3237.625
$ws.Range("K135").Value2 = 29138.625
$ws.Range("M135").Value2 = -26603.625

$ws = $wb.Worksheets.Item(1)
$ws.Range("H137").Value2 = 3800.4827
$ws.Range("I137").Value2 = 2343.2
$ws.Range("K137").Value2 = 7029.599999999999
$ws.Range("M137").Value2 = -4479.599999999999

$ws = $wb.Worksheets.Item(1)
$ws.Range("H138").Value2 = 2041.6666
$ws.Range("J138").Value2 = 3600
$ws.Range("L138").Value2 = 10800
$ws.Range("N138").Value2 = -21080

$ws = $wb.Worksheets.Item(1)
$ws.Range("H141").Value2 = 1854.7
$ws.Range("I141").Value2 = 1826.8889
$ws.Range("K141").Value2 = 5480.6667
$ws.Range("M141").Value2 = -300.6666999999998

$ws = $wb.Worksheets.Item(2)
$ws.Range("H4").Value2 = 703.6
$ws.Range("I4").Value2 = 128.33333
$ws.Range("K4").Value2 = 128.33333
$ws.Range("M4").Value2 = -12.33332999999999

$ws = $wb.Worksheets.Item(2)
$ws.Range("H5").Value2 = 510
$ws.Range("I5").Value2 = 20
$ws.Range("K5").Value2 = 20
$ws.Range("M5").Value2 = 92

$ws = $wb.Worksheets.Item(2)
$ws.Range("H51").Value2 = 0
$ws.Range("J51").Value2 = 0
$ws.Range("L51").Value2 = 0
$ws.Range("N51").ClearContents()

$ws = $wb.Worksheets.Item(2)
$ws.Range("H61").Value2 = 41756148
$ws.Range("J61").Value2 = 210849.2
$ws.Range("L61").Value2 = 210849.2
$ws.Range("N61").Value2 = -211273.2

$ws = $wb.Worksheets.Item(2)
$ws.Range("H63").Value2 = 6334.1665

$ws = $wb.Worksheets.Item(2)
$ws.Range("H66").Value2 = 6334.1665

$ws = $wb.Worksheets.Item(2)
$ws.Range("H102").Value2 = 17827.54
$ws.Range("I102").Value2 = 24530.445
$ws.Range("K102").Value2 = 24530.445
$ws.Range("M102").Value2 = -22908.445

$ws = $wb.Worksheets.Item(2)
$ws.Range("H132").Value2 = 5106.514
$ws.Range("I132").Value2 = 2785.4092
$ws.Range("K132").Value2 = 8356.2276
$ws.Range("M132").Value2 = -5826.2276

$ws = $wb.Worksheets.Item(2)
$ws.Range("H136").Value2 = 41756148
$ws.Range("J136").Value2 = 210849.2
$ws.Range("L136").Value2 = 632547.6000000001
$ws.Range("N136").Value2 = -637647.6000000001

$ws = $wb.Worksheets.Item(3)
$ws.Range("H4").Value2 = 510
$ws.Range("I4").Value2 = 20
$ws.Range("K4").Value2 = 20
$ws.Range("M4").Value2 = 95

$ws = $wb.Worksheets.Item(3)
$ws.Range("H22").Value2 = 225
$ws.Range("I22").Value2 = 225
$ws.Range("K22").Value2 = 225
$ws.Range("M22").Value2 = -52

$ws = $wb.Worksheets.Item(3)
$ws.Range("H99").Value2 = 1513.8667
$ws.Range("I99").Value2 = 1062.05
$ws.Range("J99").Value2 = 2417.5
$ws.Range("K99").Value2 = 1062.05
$ws.Range("L99").Value2 = 2417.5
$ws.Range("M99").Value2 = 435.95
$ws.Range("N99").Value2 = -5413.5

$ws = $wb.Worksheets.Item(3)
$ws.Range("H122").Value2 = 123000
$ws.Range("J122").Value2 = 123000
$ws.Range("L122").Value2 = 123000
$ws.Range("N122").Value2 = -132800

$ws = $wb.Worksheets.Item(4)
$ws.Range("H11").Value2 = 61003
$ws.Range("I11").Value2 = 5000
$ws.Range("K11").Value2 = 5000
$ws.Range("M11").Value2 = -4860

$ws = $wb.Worksheets.Item(4)
$ws.Range("H22").Value2 = 475.5
$ws.Range("J22").Value2 = 561.25
$ws.Range("L22").Value2 = 561.25
$ws.Range("N22").Value2 = -1261.25

$ws = $wb.Worksheets.Item(4)
$ws.Range("H99").Value2 = 3736.2222
$ws.Range("I99").Value2 = 3321.1667
$ws.Range("J99").Value2 = 4566.3335
$ws.Range("K99").Value2 = 3321.1667
$ws.Range("L99").Value2 = 4566.3335
$ws.Range("M99").Value2 = -1823.1667
$ws.Range("N99").Value2 = -7562.3335

$ws = $wb.Worksheets.Item(4)
$ws.Range("H126").Value2 = 3736.2222
$ws.Range("I126").Value2 = 3321.1667
$ws.Range("J126").Value2 = 4566.3335
$ws.Range("K126").Value2 = 9963.500100000001
$ws.Range("L126").Value2 = 13699.0005
$ws.Range("M126").Value2 = -7493.500100000001
$ws.Range("N126").Value2 = -18639.0005

$ws = $wb.Worksheets.Item(4)
$ws.Range("H134").Value2 = 717511.6
$ws.Range("I134").Value2 = 835012.7
$ws.Range("K134").Value2 = 2505038.1
$ws.Range("M134").Value2 = -2502503.1

$ws = $wb.Worksheets.Item(5)
$ws.Range("H2").Value2 = 187.53847
$ws.Range("I2").Value2 = 151
$ws.Range("J2").Value2 = 210.375
$ws.Range("K2").Value2 = 906
$ws.Range("L2").Value2 = 1262.25
$ws.Range("M2").Value2 = -793
$ws.Range("N2").Value2 = -1488.25

$ws = $wb.Worksheets.Item(5)
$ws.Range("H4").Value2 = 79703750
$ws.Range("J4").Value2 = 25500000
$ws.Range("L4").Value2 = 76500000
$ws.Range("N4").Value2 = -76500224

$ws = $wb.Worksheets.Item(5)
$ws.Range("H23").Value2 = 445.25
$ws.Range("I23").Value2 = 456.25
$ws.Range("J23").Value2 = 434.25
$ws.Range("K23").Value2 = 1368.75
$ws.Range("L23").Value2 = 1302.75
$ws.Range("M23").Value2 = -1133.75
$ws.Range("N23").Value2 = -1772.75

$ws = $wb.Worksheets.Item(5)
$ws.Range("H129").Value2 = 37149116
$ws.Range("I129").Value2 = 1863
$ws.Range("J129").Value2 = 55722744
$ws.Range("K129").Value2 = 5589
$ws.Range("L129").Value2 = 167168232
$ws.Range("M129").Value2 = -589
$ws.Range("N129").Value2 = -167178232

$ws = $wb.Worksheets.Item(6)
$ws.Range("H107").Value2 = 923.2353000000001
$ws.Range("J107").Value2 = 1184.7142
$ws.Range("L107").Value2 = 1184.7142
$ws.Range("N107").Value2 = -5024.7142

$ws = $wb.Worksheets.Item(6)
$ws.Range("H132").Value2 = 33336550
$ws.Range("I132").Value2 = 38464600
$ws.Range("K132").Value2 = 115393800
$ws.Range("M132").Value2 = -115391270

$ws = $wb.Worksheets.Item(7)
$ws.Range("H132").Value2 = 136185.8
$ws.Range("I132").Value2 = 93618
$ws.Range("J132").Value2 = 253247.25
$ws.Range("K132").Value2 = 280854
$ws.Range("L132").Value2 = 759741.75
$ws.Range("M132").Value2 = -278324
$ws.Range("N132").Value2 = -764801.75

$ws = $wb.Worksheets.Item(7)
$ws.Range("H136").Value2 = 197364.81
$ws.Range("J136").Value2 = 165644.28
$ws.Range("L136").Value2 = 496932.84
$ws.Range("N136").Value2 = -502032.84

$ws = $wb.Worksheets.Item(8)
$ws.Range("H126").Value2 = 4731.2
$ws.Range("I126").Value2 = 4965.8823
$ws.Range("J126").Value2 = 3401.3333
$ws.Range("K126").Value2 = 14897.6469
$ws.Range("L126").Value2 = 10203.9999
$ws.Range("M126").Value2 = -12427.6469
$ws.Range("N126").Value2 = -15143.9999

$ws = $wb.Worksheets.Item(8)
$ws.Range("H133").Value2 = 99125
$ws.Range("J133").Value2 = 99125
$ws.Range("L133").Value2 = 99125
$ws.Range("N133").Value2 = -109245

$ws = $wb.Worksheets.Item(8)
$ws.Range("H136").Value2 = 1382.4166
$ws.Range("I136").Value2 = 1067.4814
$ws.Range("J136").Value2 = 2327.2222
$ws.Range("K136").Value2 = 3202.4442
$ws.Range("L136").Value2 = 6981.6666
$ws.Range("M136").Value2 = -652.4441999999999
$ws.Range("N136").Value2 = -12081.6666
